# Assignment 2 - Task 3 - Week 10 - Answer
# Update product backlog: GUI layout manager task renamed/added, new
# "Load Images to GUI" and "Test the gameboard" tasks, progress updates.

$wb = $excel.ActiveWorkbook
$wsUserStories = $wb.Worksheets.Item("User stories")
$wsSnakeGame = $wb.Worksheets.Item("Snake-game")

# --- Snake-game sheet: product backlog updates ---

# 3.1.3 "Snake" level-3 task is renamed to "GUI layoutmanager"
$wsSnakeGame.Range("D14").Value = "GUI layoutmanager"

# New level-3 task 3.1.4 "Load Images to GUI"
$wsSnakeGame.Range("A15").Value = "3.1.4"
$wsSnakeGame.Range("D15").Value = "Load Images to GUI"

# Progress update for "Test and deploy the game" (row 25): 20% -> 30%
$wsSnakeGame.Range("F25").Value = 0.3

# Insert new sub-task row 27: "Test the gameboard" at 50% progress;
# existing row 27 ("5 - Document the game") shifts down to row 28.
$wsSnakeGame.Rows.Item(27).Insert()
$wsSnakeGame.Range("A27").Value = 4.2
$wsSnakeGame.Range("C27").Value = "Test the gameboard"
$wsSnakeGame.Range("F27").Value = 0.5

# Make the Snake-game sheet the active tab/sheet, with G16 selected.
$wsSnakeGame.Activate() | Out-Null
$wsSnakeGame.Range("G16").Select() | Out-Null
